$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 0.1079986666666667
$ws.Range("H2").Value = 0.323996
$ws.Range("I2").Value = 0.004187739561209694
$ws.Range("J2").Value = 0.004187739561209694
$ws.Range("M2").Value = 0.01393633333333333
$ws.Range("N2").Value = 0.041809
$ws.Range("Q2").Value = 0.001505105418222222
$ws.Range("R2").Value = 0.013545948764
$ws.Range("S2").Value = 0.004187739561209694
$ws.Range("T2").Value = 0.004187739561209694

# Row 3
$ws.Range("D3").Value = "MuSCs"
$ws.Range("I3").Value = 0.9687110856121154
$ws.Range("J3").Value = 0.9687110856121155
$ws.Range("M3").Value = 0.01393633333333333
$ws.Range("N3").Value = 0.041809
$ws.Range("Q3").Value = 0.3481621247777778
$ws.Range("R3").Value = 3.133459123
$ws.Range("S3").Value = 0.9687110856121154
$ws.Range("T3").Value = 0.9687110856121155

# Row 4
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6989190000000001
$ws.Range("H4").Value = 2.096757
$ws.Range("I4").Value = 0.02710117482667488
$ws.Range("J4").Value = 0.02710117482667489
$ws.Range("M4").Value = 0.01393633333333333
$ws.Range("N4").Value = 0.041809
$ws.Range("Q4").Value = 0.009740368157000001
$ws.Range("R4").Value = 0.087663313413
$ws.Range("S4").Value = 0.02710117482667488
$ws.Range("T4").Value = 0.02710117482667489
